$wb = $excel.ActiveWorkbook
$mar18 = $wb.Worksheets.Item("mar18")
$newSheet = $wb.Worksheets.Add($mar18)
$newSheet.Name = "feb18"

# --- Round 1: 3/25/2018 (serial 43184) ---
$newSheet.Range("A1").NumberFormat = "d-mmm-yy"
$newSheet.Range("A1").Value = 43184
$newSheet.Range("B1").Value = "Score"
$newSheet.Range("C1").Value = "Fairway"
$newSheet.Range("D1").Value = "GIR"
$newSheet.Range("E1").Value = "Putts"
$newSheet.Range("F1").Value = "Comment"

$newSheet.Range("A2").Value = "Hole 1"
$newSheet.Range("B2").Value = 4
$newSheet.Range("C2").Value = "R"
$newSheet.Range("E2").Value = 1

$newSheet.Range("A3").Value = "Hole 2"
$newSheet.Range("B3").Value = 4
$newSheet.Range("E3").Value = 1

$newSheet.Range("A4").Value = "Hole 3"
$newSheet.Range("B4").Value = 4
$newSheet.Range("C4").Value = "S"
$newSheet.Range("E4").Value = 2

$newSheet.Range("A5").Value = "Hole 4"
$newSheet.Range("B5").Value = 6
$newSheet.Range("C5").Value = "R"
$newSheet.Range("E5").Value = 2

$newSheet.Range("A6").Value = "Hole 5"
$newSheet.Range("B6").Value = 6
$newSheet.Range("E6").Value = 3

$newSheet.Range("A7").Value = "Hole 6"
$newSheet.Range("A8").Value = "Hole 7"
$newSheet.Range("A9").Value = "Hole 8"

$newSheet.Range("A10").Value = "Hole 9"
$newSheet.Range("B10").Value = 5
$newSheet.Range("C10").Value = "S"
$newSheet.Range("E10").Value = 2

$newSheet.Range("A11").Value = "Hole 10"
$newSheet.Range("B11").Value = 5
$newSheet.Range("C11").Value = "S"
$newSheet.Range("E11").Value = 1

$newSheet.Range("A12").Value = "Hole 11"
$newSheet.Range("B12").Value = 4
$newSheet.Range("E12").Value = 2

$newSheet.Range("A13").Value = "Hole 12"
$newSheet.Range("B13").Value = 5
$newSheet.Range("C13").Value = "R"
$newSheet.Range("E13").Value = 2

$newSheet.Range("A14").Value = "Hole 13"
$newSheet.Range("B14").Value = 5
$newSheet.Range("C14").Value = "R"
$newSheet.Range("E14").Value = 2

$newSheet.Range("A15").Value = "Hole 14"
$newSheet.Range("B15").Value = 5
$newSheet.Range("C15").Value = "L"
$newSheet.Range("E15").Value = 1

$newSheet.Range("A16").Value = "Hole 15"
$newSheet.Range("B16").Value = 3
$newSheet.Range("E16").Value = 1

$newSheet.Range("A17").Value = "Hole 16"
$newSheet.Range("B17").Value = 7
$newSheet.Range("C17").Value = "R"
$newSheet.Range("E17").Value = 2

$newSheet.Range("A18").Value = "Hole 17"
$newSheet.Range("B18").Value = 5
$newSheet.Range("C18").Value = "L"
$newSheet.Range("E18").Value = 1

$newSheet.Range("A19").Value = "Hole 18"
$newSheet.Range("B19").Value = 6
$newSheet.Range("C19").Value = "R"
$newSheet.Range("E19").Value = 2

$newSheet.Range("B20").Formula = "=SUM(B2:B19)"
$newSheet.Range("E20").Formula = "=SUM(E2:E19)"

# --- Round 2: 2/18/2018 (serial 43149) ---
$newSheet.Range("A22").NumberFormat = "d-mmm-yy"
$newSheet.Range("A22").Value = 43149
$newSheet.Range("B22").Value = "Score"
$newSheet.Range("C22").Value = "Fairway"
$newSheet.Range("D22").Value = "GIR"
$newSheet.Range("E22").Value = "Putts"
$newSheet.Range("F22").Value = "Comment"

$newSheet.Range("A23").Value = "Hole 1"
$newSheet.Range("B23").Value = 4
$newSheet.Range("C23").Value = "R"
$newSheet.Range("E23").Value = 1

$newSheet.Range("A24").Value = "Hole 2"
$newSheet.Range("B24").Value = 3
$newSheet.Range("E24").Value = 2

$newSheet.Range("A25").Value = "Hole 3"
$newSheet.Range("B25").Value = 4
$newSheet.Range("C25").Value = "S"
$newSheet.Range("E25").Value = 1

$newSheet.Range("A26").Value = "Hole 4"
$newSheet.Range("B26").Value = 5
$newSheet.Range("C26").Value = "S"
$newSheet.Range("E26").Value = 1

$newSheet.Range("A27").Value = "Hole 5"
$newSheet.Range("B27").Value = 3
$newSheet.Range("E27").Value = 2

$newSheet.Range("A28").Value = "Hole 6"
$newSheet.Range("B28").Value = 6
$newSheet.Range("C28").Value = "S"
$newSheet.Range("E28").Value = 2

$newSheet.Range("A29").Value = "Hole 7"
$newSheet.Range("B29").Value = 5
$newSheet.Range("C29").Value = "R"
$newSheet.Range("E29").Value = 2

$newSheet.Range("A30").Value = "Hole 8"
$newSheet.Range("B30").Value = 3
$newSheet.Range("E30").Value = 2

$newSheet.Range("A31").Value = "Hole 9"
$newSheet.Range("B31").Value = 5
$newSheet.Range("C31").Value = "R"
$newSheet.Range("E31").Value = 2

$newSheet.Range("A32").Value = "Hole 10"
$newSheet.Range("B32").Value = 4
$newSheet.Range("C32").Value = "R"
$newSheet.Range("E32").Value = 1

$newSheet.Range("A33").Value = "Hole 11"
$newSheet.Range("B33").Value = 4
$newSheet.Range("E33").Value = 2

$newSheet.Range("A34").Value = "Hole 12"
$newSheet.Range("B34").Value = 5
$newSheet.Range("C34").Value = "S"
$newSheet.Range("E34").Value = 2

$newSheet.Range("A35").Value = "Hole 13"
$newSheet.Range("A36").Value = "Hole 14"

$newSheet.Range("A37").Value = "Hole 15"
$newSheet.Range("B37").Value = 3
$newSheet.Range("E37").Value = 2

$newSheet.Range("A38").Value = "Hole 16"
$newSheet.Range("A39").Value = "Hole 17"

$newSheet.Range("A40").Value = "Hole 18"
$newSheet.Range("B40").Value = 5
$newSheet.Range("C40").Value = "R"
$newSheet.Range("E40").Value = 2

$newSheet.Range("B41").Formula = "=SUM(B23:B40)"
$newSheet.Range("E41").Formula = "=SUM(E23:E40)"

$newSheet.Range("A43").NumberFormat = "d-mmm-yy"

# --- View settings to match target (feb18 becomes the active/selected tab, selection B40) ---
$newSheet.Range("B40").Select()
